$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume (E) columns keep their text representation
# (values like "0.530" or "174.60" would otherwise be coerced to numbers)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "71.469.81"
$ws.Range("E2").Value = "  +1.15%  "

$ws.Range("D3").Value = "3.813.03"
$ws.Range("E3").Value = "  -0.04%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "699.95"
$ws.Range("E5").Value = "  +5.32%  "

$ws.Range("D6").Value = "174.60"
$ws.Range("E6").Value = "  +3.38%  "

$ws.Range("D7").Value = "3.812.04"
$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").Value = "0.530"
$ws.Range("E9").Value = "  +0.32%  "

$ws.Range("E10").Value = "  +1.31%  "

$ws.Range("D11").Value = "7.43"
$ws.Range("E11").Value = "  +6.15%  "

$ws.Range("E12").Value = "  +0.08%  "

$ws.Range("D13").Value = "0.0000259"
$ws.Range("E13").Value = "  +6.11%  "

$ws.Range("D14").Value = "36.67"
$ws.Range("E14").Value = "  +2.20%  "

$ws.Range("D15").Value = "4.443.60"
$ws.Range("E15").Value = "  -0.28%  "

$ws.Range("D16").Value = "3.795.88"
$ws.Range("E16").Value = "  -0.55%  "

$ws.Range("D17").Value = "71.382.63"
$ws.Range("E17").Value = "  +1.08%  "

$ws.Range("D18").Value = "17.78"
$ws.Range("E18").Value = "  +0.03%  "

$ws.Range("D19").Value = "7.25"
$ws.Range("E19").Value = "  +1.21%  "

$ws.Range("E20").Value = "  +0.45%  "

$ws.Range("D21").Value = "11.14"
$ws.Range("E21").Value = "  +7.53%  "

$ws.Range("D22").Value = "487.40"
$ws.Range("E22").Value = "  +1.79%  "

$ws.Range("D23").Value = "0.718"
$ws.Range("E23").Value = "  +0.71%  "

$ws.Range("D24").Value = "84.78"
$ws.Range("E24").Value = "  +2.34%  "

$ws.Range("D25").Value = "0.0000143"
$ws.Range("E25").Value = "  -1.42%  "

$ws.Range("D26").Value = "12.36"
$ws.Range("E26").Value = "  +0.70%  "

$ws.Range("D27").Value = "10.56"
$ws.Range("E27").Value = "  +1.68%  "

$ws.Range("E28").Value = "  +1.71%  "

$ws.Range("D29").Value = "3.960.87"
$ws.Range("E29").Value = "  -0.13%  "

$ws.Range("E30").Value = "  +0.11%  "

$ws.Range("D31").Value = "3.13"

$ws.Range("D32").Value = "2.32"
$ws.Range("E32").Value = "  +0.43%  "

$ws.Range("D33").Value = "7.62"
$ws.Range("E33").Value = "  +2.92%  "

$ws.Range("D34").Value = "0.185"
$ws.Range("E34").Value = "  +2.72%  "

$ws.Range("D35").Value = "29.69"
$ws.Range("E35").Value = "  +0.49%  "

$ws.Range("D36").Value = "9.31"
$ws.Range("E36").Value = "  +1.76%  "

$ws.Range("E37").Value = "  +0.06%  "

$ws.Range("D38").Value = "0.104"
$ws.Range("E38").Value = "  +2.03%  "

$ws.Range("D39").Value = "2.39"
$ws.Range("E39").Value = "  +15.21%  "

$ws.Range("E40").Value = "  -0.18%  "

$ws.Range("D41").Value = "6.05"
$ws.Range("E41").Value = "  +2.21%  "

$ws.Range("D42").Value = "0.998"
$ws.Range("E42").Value = "  +2.60%  "

$ws.Range("D43").Value = "0.998"
$ws.Range("E43").Value = "  -0.31%  "

$ws.Range("E44").Value = "  +0.03%  "

$ws.Range("D45").Value = "164.79"
$ws.Range("E45").Value = "  +3.91%  "

$ws.Range("D46").Value = "0.000306"
$ws.Range("E46").Value = "  +6.75%  "

$ws.Range("D47").Value = "44.69"
$ws.Range("E47").Value = "  -1.65%  "

$ws.Range("D48").Value = "48.66"
$ws.Range("E48").Value = "  -0.23%  "

$ws.Range("B49").Value = "Bittensor"
$ws.Range("C49").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D49").Value = "419.05"
$ws.Range("E49").Value = "  +5.93%  "

$ws.Range("B50").Value = "TheGraph"
$ws.Range("C50").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D50").Value = "0.303"
$ws.Range("E50").Value = "  +1.13%  "

$ws.Range("D51").Value = "8.69"
$ws.Range("E51").Value = "  +2.29%  "
